# Today_words.xlsx - "add new 15 words from office belong to D cat"
# Adds 15 English words (column B) with their Hindi translations (column D)
# to rows 1-15 of Sheet1, resizes columns C/D, clears a stray styled cell,
# and updates the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (word) / Column D (Hindi meaning) pairs, row by row -------
# Values are written in the same B/D interleave order that produced the
# original shared-string table, so the resulting sharedStrings.xml lines
# up with the source edit as closely as possible.

$ws.Range("B1").Value = "goodness sake"
$ws.Range("D1").Value = "भलाई के"

$ws.Range("B2").Value = "deceased"
$ws.Range("D2").Value = "मृतक"

$ws.Range("D3").Value = "रोगी"
$ws.Range("B3").Value = "diseased"

$ws.Range("D4").Value = "छल"
$ws.Range("B4").Value = "deceit"

$ws.Range("D5").Value = "धोखा देना"
$ws.Range("B5").Value = "deceive"

$ws.Range("D6").Value = "सभ्य"
$ws.Range("B6").Value = "decent"

$ws.Range("B7").Value = "descent"
$ws.Range("D7").Value = "अवरोह"

$ws.Range("D8").Value = "संकेतमय"
$ws.Range("B8").Value = "ALLUSION"

$ws.Range("B9").Value = "DELUSION"
$ws.Range("D9").Value = "भ्रम"

$ws.Range("D10").Value = "सूखा"
$ws.Range("B10").Value = "desiccated"

$ws.Range("D11").Value = "निराश"
$ws.Range("B11").Value = "desperate"

$ws.Range("D12").Value = "अलग"
$ws.Range("B12").Value = "detached"

$ws.Range("D13").Value = "आपदा"
$ws.Range("B13").Value = "disaster"

$ws.Range("B14").Value = "disasterous"
$ws.Range("D14").Value = "विनाशकारी"

$ws.Range("D15").Value = "विचारशील"
$ws.Range("B15").Value = "discreet"

# --- Formatting ------------------------------------------------------
# Column B keeps the existing blue "word" style (font color 1F497D).
# Only touch Color so the cell reuses the workbook's existing font entry
# instead of minting a near-duplicate (losing the Calibri "minor" scheme).
$ws.Range("B1:B15").Font.Color = 8210719

# D7 shares the small black Arial "meaning" style already used by the
# C column helper cells (C1, C2, C7, C11, C16).
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 10
$ws.Range("D7").Font.Color = 0

# Column widths: C narrows, D is a brand-new column.
$ws.Columns.Item(3).ColumnWidth = 12.6
$ws.Columns.Item(4).ColumnWidth = 15.6

# The old C49 placeholder cell is no longer part of the filled range.
$ws.Range("C49").Clear()

# Leave the view scrolled near the new data, with B15 selected, matching
# the author's final cursor position.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()
